$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeviceList")

# Remove the device columns for devices no longer part of this rerun batch:
# SAMSUNG_GalaxyFold (col C), SAMSUNG_GalaxyNote20_12.0.0 (col E),
# SAMSUNG_GalaxyTabA7 (col H), SAMSUNG_GalaxyA32 (col I)
# Deleting right-to-left keeps the remaining column letters stable while we work.
$ws.Columns("I").Delete()
$ws.Columns("H").Delete()
$ws.Columns("E").Delete()
$ws.Columns("C").Delete()

# Move the active selection to match where the author left off editing.
$ws.Range("C17").Select()
